# Commit funcional de ejecucion bd y soap
# Applies the data/structure changes described by the diff:
#  - Removes the "URL GATEWAY MG" column (E) and its hyperlink
#  - Fixes header D10 ("MSI" -> "IMSI")
#  - Fixes B12 / B14 values (ICCID had been truncated to the IMSI value)
#  - Adds new data rows 15-18 (additional SIM/IMSI test records)
#  - Removes the stray empty styled cell D19
#  - Rebuilds hyperlinks (D2, A2) and selection/view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove hyperlinks so we can rebuild only the ones that remain
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Drop the "URL GATEWAY MG" column (E) entirely; columns F:J were already
#    empty helper cells beyond the real table, so remove them too.
# ---------------------------------------------------------------------------
$ws.Columns("E:J").Delete()

# ---------------------------------------------------------------------------
# 3. Fix header + data values that changed in place
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = "IMSI"
$ws.Range("B12").Value = "8957732111198172293"
$ws.Range("B14").Value = "8957732111198172293"

# ---------------------------------------------------------------------------
# 4. Helper to write the PLU number so it stays a real number (not text)
#    even though the column is formatted as Text (numFmtId 49 / "@"),
#    matching the source file's pre-existing convention for column A.
# ---------------------------------------------------------------------------
function Set-NumericTextCell($cell, $number) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "General"
    $cell.Value2 = $number
    $cell.NumberFormat = $fmt
}

# ---------------------------------------------------------------------------
# 5. Add the new records (rows 15-18)
# ---------------------------------------------------------------------------
Set-NumericTextCell $ws.Range("A15") 3003324
$ws.Range("B15").Value = "8957732111198172291"
$ws.Range("C15").Value = "3016877412"
$ws.Range("D15").Value = "732111198172299"

Set-NumericTextCell $ws.Range("A16") 3003324
$ws.Range("B16").Value = "8957732111198172290"
$ws.Range("C16").Value = "3016876873"
$ws.Range("D16").Value = "732111198172297"

Set-NumericTextCell $ws.Range("A17") 3003324
$ws.Range("B17").Value = "8957732111198172295"
$ws.Range("C17").Value = "3016877414"
$ws.Range("D17").Value = "732111198172296"

Set-NumericTextCell $ws.Range("A18") 3003324
$ws.Range("B18").Value = "8957732111198172296"
$ws.Range("C18").Value = "3016876878"
$ws.Range("D18").Value = "732111198172295"

# ---------------------------------------------------------------------------
# 6. Remove the stray formatted-but-empty cell that used to sit at D19
# ---------------------------------------------------------------------------
$ws.Rows("19:19").Delete()

# ---------------------------------------------------------------------------
# 7. Rebuild the two remaining hyperlinks (gateway CBS + EPOS) in order.
#    Adding a hyperlink re-applies Excel's built-in hyperlink look, so we
#    restore the original cell formatting afterwards by copying the format
#    from untouched cells that still carry the same original style.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D2"), "http://10.65.45.12:9001/gatewaycbs/BcServicesInt")
$ws.Hyperlinks.Add($ws.Range("A2"), "http://10.69.60.77:8180/tigo-pos-web/")

$ws.Range("C8").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 8. Restore view state: scrolled so row 4 is at the top, selection on E19
# ---------------------------------------------------------------------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E19").Select()
